# param_data_naninovel.xlsx update
# - Pose sheet: new row for "姿势变换" / "the_pose"
# - Tint sheet: new row for "色调" / "the_tint_color"
# - PackedEffect sheet: updated ShakeCamera / Blur command strings
# - PackedEffect sheet becomes the active/selected sheet (was Printer)

$wb = $excel.ActiveWorkbook

# --- Pose: insert "姿势变换 / the_pose" as the new second row -----------------
$wsPose = $wb.Worksheets.Item("Pose")
$wsPose.Rows.Item(2).Insert()
$wsPose.Range("A2").Value = "姿势变换"
$wsPose.Range("B2").Value = "the_pose"
$wsPose.Range("D17").Select()

# --- Tint: insert "色调 / the_tint_color" as the new second row --------------
$wsTint = $wb.Worksheets.Item("Tint")
$wsTint.Rows.Item(2).Insert()
$wsTint.Range("A2").Value = "色调"
$wsTint.Range("B2").Value = "the_tint_color"

# --- PackedEffect: refresh the generated command strings ---------------------
$wsPacked = $wb.Worksheets.Item("PackedEffect")
$wsPacked.Range("B2").Value = "'@spawn ShakeCamera params:,2,,,0.3,"
$wsPacked.Range("B3").Value = "'@spawn ShakeCamera params:,3,,,0.6,"
$wsPacked.Range("B4").Value = "'@blur power:0"
$wsPacked.Range("B5").Select()

# --- Printer: move selection, no longer the active tab -----------------------
$wsPrinter = $wb.Worksheets.Item("Printer")
$wsPrinter.Range("B10").Select()

# --- PackedEffect becomes the active sheet/tab -------------------------------
$wsPacked.Activate()
$wsPacked.Range("B5").Select()
